# Apply the "dashboard completo" edit:
#  - Insert two new blank columns before the old "Demanda diaria" column
#    (old M:Q shift right to become O:S)
#  - Label the two new columns "Unnamed: 12" / "Unnamed: 13"
#  - Clear out the existing "Ventas" values in column K (header stays)
#  - The new M column stays blank; the new N column only gets one value,
#    in row 2 (2840000) - all the other rows stay blank in N

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns starting at M (shifts old M..Q to O..S)
$ws.Range("M1:N1").EntireColumn.Insert()

# The insert copies neighboring column formatting (e.g. the date style from
# column L) onto the new cells - strip that back off so the new columns'
# data cells are unformatted like the rest of the plain data columns.
$ws.Range("M2:N14").ClearFormats()

# Header labels for the two newly inserted columns
$ws.Range("M1").Value = "Unnamed: 12"
$ws.Range("N1").Value = "Unnamed: 13"

# Clear the old "Ventas" data values (column K) while keeping its header
$ws.Range("K2:K14").ClearContents()

# New data point introduced only for row 2
$ws.Range("N2").Value = 2840000
